# Apply the refreshed crypto price/volume figures (and the Stellar/Arweave row swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.612.97'
$ws.Range("E2").Value = '  +0.83%  '
$ws.Range("D3").Value = '3.291.55'
$ws.Range("E3").Value = '  +4.74%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '''600.55'
$ws.Range("E5").Value = '  +2.50%  '
$ws.Range("D6").Value = '''141.74'
$ws.Range("E6").Value = '  +3.20%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = '3.290.40'
$ws.Range("E8").Value = '  +4.88%  '
$ws.Range("E9").Value = '  +0.59%  '
$ws.Range("E10").Value = '  +2.96%  '
$ws.Range("D11").Value = '''5.44'
$ws.Range("E11").Value = '  +3.84%  '
$ws.Range("D12").Value = '''0.469'
$ws.Range("E12").Value = '  +2.57%  '
$ws.Range("D13").Value = '''0.0000247'
$ws.Range("E13").Value = '  +1.02%  '
$ws.Range("D14").Value = '''34.49'
$ws.Range("E14").Value = '  +0.94%  '
$ws.Range("D15").Value = '3.836.10'
$ws.Range("E15").Value = '  +4.84%  '
$ws.Range("E16").Value = '  +0.96%  '
$ws.Range("D17").Value = '3.288.75'
$ws.Range("E17").Value = '  +4.84%  '
$ws.Range("D18").Value = '63.649.16'
$ws.Range("E18").Value = '  +0.98%  '
$ws.Range("D19").Value = '''6.83'
$ws.Range("E19").Value = '  +2.70%  '
$ws.Range("D20").Value = '''477.90'
$ws.Range("E20").Value = '  +1.35%  '
$ws.Range("D21").Value = '''14.08'
$ws.Range("E21").Value = '  -0.42%  '
$ws.Range("E22").Value = '  +4.40%  '
$ws.Range("D23").Value = '''8.05'
$ws.Range("E23").Value = '  +5.09%  '
$ws.Range("D24").Value = '''13.61'
$ws.Range("E24").Value = '  +5.01%  '
$ws.Range("D25").Value = '''84.18'
$ws.Range("E25").Value = '  -0.66%  '
$ws.Range("E26").Value = '  +0.08%  '
$ws.Range("E27").Value = '  +2.06%  '
$ws.Range("E28").Value = '  -0.23%  '
$ws.Range("D29").Value = '''7.22'
$ws.Range("E29").Value = '  +3.25%  '
$ws.Range("D30").Value = '''8.08'
$ws.Range("E30").Value = '  +1.71%  '
$ws.Range("E31").Value = '  +1.75%  '
$ws.Range("D32").Value = '''28.61'
$ws.Range("E32").Value = '  +6.99%  '
$ws.Range("E33").Value = '  -1.42%  '
$ws.Range("E34").Value = '  +0.30%  '
$ws.Range("E35").Value = '  +3.42%  '
$ws.Range("D36").Value = '''5.97'
$ws.Range("E36").Value = '  +3.45%  '
$ws.Range("D37").Value = '''53.17'
$ws.Range("E37").Value = '  +1.78%  '
$ws.Range("E38").Value = '  +6.45%  '
$ws.Range("E39").Value = '  +3.03%  '
$ws.Range("D40").Value = '''427.13'
$ws.Range("E40").Value = '  +2.15%  '
$ws.Range("D41").Value = '3.071.00'
$ws.Range("E41").Value = '  +6.14%  '
$ws.Range("D42").Value = '''8.33'
$ws.Range("E42").Value = '  +1.82%  '
$ws.Range("E43").Value = '  +0.17%  '
$ws.Range("E44").Value = '  +0.90%  '
$ws.Range("D45").Value = '''0.264'
$ws.Range("E45").Value = '  +1.25%  '
$ws.Range("E46").Value = '  +3.37%  '
$ws.Range("D47").Value = '''26.21'
$ws.Range("E47").Value = '  +3.15%  '
$ws.Range("D49").Value = '''127.08'
$ws.Range("E49").Value = '  +5.54%  '
$ws.Range("B50").Value = 'Arweave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D50").Value = '''35.47'
$ws.Range("E50").Value = '  +13.32%  '
$ws.Range("B51").Value = 'Stellar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D51").Value = '''0.114'
$ws.Range("E51").Value = '  +1.54%  '
